$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1648.4166
$ws.Range("I12").Value = 1630.1111
$ws.Range("J12").Value = 1703.3334
$ws.Range("K12").Value = 1630.1111
$ws.Range("L12").Value = 1703.3334
$ws.Range("M12").Value = -1460.1111
$ws.Range("N12").Value = -2043.3334

$ws.Range("H32").Value = 5599.3
$ws.Range("J32").Value = 6499.25
$ws.Range("L32").Value = 6499.25
$ws.Range("N32").Value = -7151.25

$ws.Range("H33").Value = 11365798
$ws.Range("I33").Value = 15626244
$ws.Range("K33").Value = 15626244
$ws.Range("M33").Value = -15626015

$ws.Range("H43").Value = 1787.26
$ws.Range("I43").Value = 1497.034
$ws.Range("J43").Value = 3915.5833
$ws.Range("K43").Value = 1497.034
$ws.Range("L43").Value = 3915.5833
$ws.Range("M43").Value = -1428.034
$ws.Range("N43").Value = -4053.5833

$ws.Range("H112").Value = 1989.5918
$ws.Range("J112").Value = 2088.6667
$ws.Range("L112").Value = 6266.000100000001
$ws.Range("N112").Value = -8482.000100000001

$ws.Range("H113").Value = 16993.908
$ws.Range("I113").Value = 17539.666
$ws.Range("J113").Value = 16339
$ws.Range("K113").Value = 17539.666
$ws.Range("L113").Value = 16339
$ws.Range("M113").Value = -14285.666
$ws.Range("N113").Value = -22847

$ws.Range("H115").Value = 652.8
$ws.Range("J115").Value = 990
$ws.Range("L115").Value = 2970
$ws.Range("N115").Value = -6104

$ws.Range("H120").Value = 101488.836
$ws.Range("I120").Value = 100000
$ws.Range("J120").Value = 102233.25
$ws.Range("K120").Value = 100000
$ws.Range("L120").Value = 102233.25
$ws.Range("M120").Value = -95162
$ws.Range("N120").Value = -111909.25

$ws.Range("H125").Value = 3522.6667
$ws.Range("I125").Value = 3531.5715
$ws.Range("K125").Value = 31784.1435
$ws.Range("M125").Value = -29324.1435

$ws.Range("H132").Value = 15595.19
$ws.Range("I132").Value = 12964.294
$ws.Range("K132").Value = 38892.882
$ws.Range("M132").Value = -36362.882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 268.73685
$ws.Range("I4").Value = 276.88235
$ws.Range("K4").Value = 276.88235
$ws.Range("M4").Value = -160.88235

$ws.Range("H9").Value = 50000
$ws.Range("J9").Value = 50000
$ws.Range("L9").Value = 50000
$ws.Range("N9").Value = -50340

$ws.Range("H20").Value = 50000
$ws.Range("J20").Value = 50000
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50540

$ws.Range("H32").Value = 5566.1465
$ws.Range("I32").Value = 3218.224
$ws.Range("K32").Value = 3218.224
$ws.Range("M32").Value = -2931.224

$ws.Range("H132").Value = 17863.459
$ws.Range("I132").Value = 6129.3076
$ws.Range("J132").Value = 31731.092
$ws.Range("K132").Value = 18387.9228
$ws.Range("L132").Value = 95193.276
$ws.Range("M132").Value = -15857.9228
$ws.Range("N132").Value = -100253.276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 21705.426
$ws.Range("I20").Value = 6267.5264
$ws.Range("K20").Value = 6267.5264
$ws.Range("M20").Value = -6020.5264

$ws.Range("H22").Value = 1592.1428
$ws.Range("I22").Value = 2093
$ws.Range("K22").Value = 2093
$ws.Range("M22").Value = -1920

$ws.Range("H95").Value = 22549.4
$ws.Range("J95").Value = 22549.4
$ws.Range("L95").Value = 22549.4
$ws.Range("N95").Value = -28041.4

$ws.Range("H107").Value = 2632.6428
$ws.Range("I107").Value = 2895.7
$ws.Range("J107").Value = 1975
$ws.Range("K107").Value = 2895.7
$ws.Range("L107").Value = 1975
$ws.Range("M107").Value = -975.6999999999998
$ws.Range("N107").Value = -5815

$ws.Range("H132").Value = 88386.75
$ws.Range("J132").Value = 88386.75
$ws.Range("L132").Value = 88386.75
$ws.Range("N132").Value = -98506.75

$ws.Range("H134").Value = 6264.978
$ws.Range("I134").Value = 1403.8948
$ws.Range("K134").Value = 4211.6844
$ws.Range("M134").Value = -1676.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1866.4667
$ws.Range("I22").Value = 1580
$ws.Range("K22").Value = 1580
$ws.Range("M22").Value = -1230

$ws.Range("H62").Value = 3753.2
$ws.Range("J62").Value = 4066.5
$ws.Range("L62").Value = 4066.5
$ws.Range("N62").Value = -5314.5

$ws.Range("H65").Value = 3753.2
$ws.Range("J65").Value = 4066.5
$ws.Range("L65").Value = 20332.5
$ws.Range("N65").Value = -26572.5

$ws.Range("H86").Value = 24049.111
$ws.Range("I86").Value = 35597.8
$ws.Range("J86").Value = 9613.25
$ws.Range("K86").Value = 35597.8
$ws.Range("L86").Value = 9613.25
$ws.Range("M86").Value = -34474.8
$ws.Range("N86").Value = -11859.25

$ws.Range("H89").Value = 24049.111
$ws.Range("I89").Value = 35597.8
$ws.Range("J89").Value = 9613.25
$ws.Range("K89").Value = 177989
$ws.Range("L89").Value = 48066.25
$ws.Range("M89").Value = -172373
$ws.Range("N89").Value = -59298.25

$ws.Range("H134").Value = 3794.8628
$ws.Range("I134").Value = 890.4737
$ws.Range("K134").Value = 2671.4211
$ws.Range("M134").Value = -136.4211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 335.53845
$ws.Range("I46").Value = 335.53845
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1006.61535
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -915.61535
$ws.Range("N46").ClearContents()

$ws.Range("H63").Value = 15505
$ws.Range("I63").Value = 112
$ws.Range("K63").Value = 336
$ws.Range("M63").Value = 413

$ws.Range("H66").Value = 15505
$ws.Range("I66").Value = 112
$ws.Range("K66").Value = 1008
$ws.Range("M66").Value = 2736

$ws.Range("H137").Value = 11334.125
$ws.Range("J137").Value = 14130.5
$ws.Range("L137").Value = 42391.5
$ws.Range("N137").Value = -52591.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12288
$ws.Range("I43").Value = 12504
$ws.Range("K43").Value = 12504
$ws.Range("M43").Value = -12353

$ws.Range("H46").Value = 16679.166
$ws.Range("J46").Value = 34190
$ws.Range("L46").Value = 34190
$ws.Range("N46").Value = -34502

$ws.Range("H70").Value = 11959.228
$ws.Range("I70").Value = 12532.214
$ws.Range("J70").Value = 10956.5
$ws.Range("K70").Value = 12532.214
$ws.Range("L70").Value = 10956.5
$ws.Range("M70").Value = -12262.214
$ws.Range("N70").Value = -11496.5

$ws.Range("H73").Value = 11959.228
$ws.Range("I73").Value = 12532.214
$ws.Range("J73").Value = 10956.5
$ws.Range("K73").Value = 12532.214
$ws.Range("L73").Value = 10956.5
$ws.Range("M73").Value = -11596.214
$ws.Range("N73").Value = -12828.5

$ws.Range("H102").Value = 5592.263
$ws.Range("I102").Value = 6303.3125
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 6303.3125
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = -4681.3125
$ws.Range("N102").Value = -5044

$ws.Range("H122").Value = 4007.2354
$ws.Range("I122").Value = 4108.857
$ws.Range("J122").Value = 3533
$ws.Range("K122").Value = 12326.571
$ws.Range("L122").Value = 10599
$ws.Range("M122").Value = -9876.571
$ws.Range("N122").Value = -15499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2493.7
$ws.Range("I16").Value = 2309.7058
$ws.Range("K16").Value = 2309.7058
$ws.Range("M16").Value = -2139.7058

$ws.Range("H46").Value = 3602.2144
$ws.Range("I46").Value = 3187.6
$ws.Range("K46").Value = 3187.6
$ws.Range("M46").Value = -2999.6

$ws.Range("H61").Value = 3767.75
$ws.Range("I61").Value = 3252.8572
$ws.Range("J61").Value = 5312.4287
$ws.Range("K61").Value = 3252.8572
$ws.Range("L61").Value = 5312.4287
$ws.Range("M61").Value = -3050.8572
$ws.Range("N61").Value = -5716.4287

$ws.Range("H100").Value = 3728.7856
$ws.Range("I100").Value = 3677.1538
$ws.Range("K100").Value = 3677.1538
$ws.Range("M100").Value = -3136.1538

$ws.Range("H113").Value = 3767.75
$ws.Range("I113").Value = 3252.8572
$ws.Range("J113").Value = 5312.4287
$ws.Range("K113").Value = 3252.8572
$ws.Range("L113").Value = 5312.4287
$ws.Range("M113").Value = -1082.8572
$ws.Range("N113").Value = -9652.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 30226.8
$ws.Range("I62").Value = 8711.333000000001
$ws.Range("K62").Value = 8711.333000000001
$ws.Range("M62").Value = -8087.333000000001

$ws.Range("H64").Value = 111615.305
$ws.Range("J64").Value = 111615.305
$ws.Range("L64").Value = 111615.305
$ws.Range("N64").Value = -112111.305

$ws.Range("H65").Value = 30226.8
$ws.Range("I65").Value = 8711.333000000001
$ws.Range("K65").Value = 43556.665
$ws.Range("M65").Value = -40436.665

$ws.Range("H67").Value = 111615.305
$ws.Range("J67").Value = 111615.305
$ws.Range("L67").Value = 111615.305
$ws.Range("N67").Value = -113331.305

$ws.Range("H122").Value = 10692.154
$ws.Range("I122").Value = 2666.6667
$ws.Range("K122").Value = 8000.000100000001
$ws.Range("M122").Value = -5550.000100000001

$ws.Range("H132").Value = 4570.463
$ws.Range("I132").Value = 1830.4546
$ws.Range("K132").Value = 5491.3638
$ws.Range("M132").Value = -2961.3638
